$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 3 days.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Value()
    $cell.Value = $cur.AddDays(3)
}

# Updated forecast values for column B (rows 21-47), per retrained model.
$newB = @{
    21 = 6
    22 = 23
    23 = 49
    24 = 87
    25 = 129
    26 = 199
    27 = 262
    28 = 306
    29 = 393
    30 = 498
    31 = 589
    32 = 660
    33 = 783
    34 = 912
    35 = 1006
    36 = 1052
    37 = 1134
    38 = 1251
    39 = 1308
    40 = 1344
    41 = 1351
    42 = 1424
    43 = 1482
    44 = 1545
    45 = 1557
    46 = 1643
    47 = 1678
}

foreach ($r in $newB.Keys) {
    $ws.Cells.Item($r, 2).Value = $newB[$r]
}
